$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.51129999999999
$ws.Range("E3").Value = 16.3604
$ws.Range("A12").Value = -21.57660000000001
$ws.Range("C14").Value = -13.5527
$ws.Range("E20").Value = 15.95419999999999
$ws.Range("E25").Value = 17.0813
$ws.Range("C26").Value = -11.7321
$ws.Range("A27").Value = -21.80449999999999
$ws.Range("E30").Value = 15.64060000000001
$ws.Range("C31").Value = -12.5651
$ws.Range("A32").Value = -21.11699999999998
$ws.Range("C35").Value = -13.08860000000002
$ws.Range("A36").Value = -19.7436
$ws.Range("C37").Value = -13.8469
$ws.Range("A38").Value = -19.2315
$ws.Range("E44").Value = 16.75199999999999
$ws.Range("C45").Value = -13.7504
$ws.Range("A46").Value = -21.3096
$ws.Range("E47").Value = 16.42989999999999
$ws.Range("C52").Value = -10.8221
$ws.Range("A54").Value = -21.58099999999999
$ws.Range("A55").Value = -22.53220000000001
$ws.Range("A56").Value = -22.38780000000001
$ws.Range("C57").Value = -14.41869999999999
$ws.Range("E58").Value = 16.6804
$ws.Range("A67").Value = -21.55929999999999
$ws.Range("A69").Value = -21.60259999999999
$ws.Range("A72").Value = -21.37549999999999
$ws.Range("E78").Value = 16.63420000000002
$ws.Range("C81").Value = -13.06140000000001
$ws.Range("A83").Value = -21.88979999999999
$ws.Range("C83").Value = -12.1764
$ws.Range("E84").Value = 16.66500000000001
$ws.Range("A86").Value = -22.24580000000001
$ws.Range("E89").Value = 17.43910000000001
$ws.Range("A91").Value = -21.46289999999999
$ws.Range("E91").Value = 17.90380000000002
$ws.Range("E92").Value = 17.95330000000002
$ws.Range("A93").Value = -21.23599999999999
$ws.Range("E96").Value = 16.04369999999999
$ws.Range("A99").Value = -20.07969999999999
$ws.Range("C100").Value = -12.99709999999999
$ws.Range("C102").Value = -12.302
$ws.Range("E102").Value = 16.47209999999999
